$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 3 and 4: "nie" -> "tak" in column A
$ws.Range("A3").Value = "tak"
$ws.Range("A4").Value = "tak"

# New row 5
$ws.Range("A5").Value = "nie"
$ws.Range("B5").Value = 8
$ws.Range("C5").Value = "dddd"
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = "16/11/2023"
$ws.Range("F5").Value = "dd"
$ws.Range("G5").Value = ""

# New row 6
$ws.Range("A6").Value = "nie"
$ws.Range("B6").Value = 9
$ws.Range("C6").Value = "asdada"
$ws.Range("D6").Value = ""
$ws.Range("E6").Value = "16/11/2023"
$ws.Range("F6").Value = ""
$ws.Range("G6").Value = ""
